$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'47.503.76"
$ws.Range("E2").Value = '  +5.29%  '

$ws.Range("D3").Value = "'2.498.21"
$ws.Range("E3").Value = '  +2.66%  '

$ws.Range("D5").Value = "'323.52"
$ws.Range("E5").Value = '  +2.12%  '

$ws.Range("D6").Value = "'105.28"
$ws.Range("E6").Value = '  +2.22%  '

$ws.Range("D7").Value = "'0.522"
$ws.Range("E7").Value = '  +1.47%  '

$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("D9").Value = "'0.539"
$ws.Range("E9").Value = '  +2.42%  '

$ws.Range("D10").Value = "'37.64"
$ws.Range("E10").Value = '  +6.08%  '

$ws.Range("E11").Value = '  +1.32%  '

$ws.Range("E12").Value = '  +0.27%  '

$ws.Range("E13").Value = '  +0.47%  '

$ws.Range("D14").Value = "'7.15"
$ws.Range("E14").Value = '  +2.13%  '

$ws.Range("D15").Value = "'2.885.23"
$ws.Range("E15").Value = '  +2.71%  '

$ws.Range("D16").Value = "'2.496.53"
$ws.Range("E16").Value = '  +3.26%  '

$ws.Range("E17").Value = '  +0.57%  '

$ws.Range("D18").Value = "'47.391.26"
$ws.Range("E18").Value = '  +5.27%  '

$ws.Range("E19").Value = '  +3.51%  '

$ws.Range("E20").Value = '  +2.71%  '

$ws.Range("D21").Value = "'0.0₃0935"
$ws.Range("E21").Value = '  +1.20%  '

$ws.Range("D22").Value = "'70.85"
$ws.Range("E22").Value = '  +2.88%  '

$ws.Range("D23").Value = "'250.75"
$ws.Range("E23").Value = '  +2.86%  '

$ws.Range("D24").Value = "'2.40"
$ws.Range("E24").Value = '  +5.79%  '

$ws.Range("D25").Value = "'2.57"
$ws.Range("E25").Value = '  +2.90%  '

$ws.Range("D26").Value = "'26.19"
$ws.Range("E26").Value = '  +3.57%  '

$ws.Range("E27").Value = '  -0.03%  '

$ws.Range("E28").Value = '  +5.47%  '

$ws.Range("D29").Value = "'2.22"
$ws.Range("E29").Value = '  +1.50%  '

$ws.Range("D30").Value = "'35.25"
$ws.Range("E30").Value = '  +7.21%  '

$ws.Range("D31").Value = "'0.134"
$ws.Range("E31").Value = '  +8.66%  '

$ws.Range("D32").Value = "'49.48"
$ws.Range("E32").Value = '  +0.48%  '

$ws.Range("D33").Value = "'20.00"
$ws.Range("E33").Value = '  -0.73%  '

$ws.Range("D34").Value = "'5.37"
$ws.Range("E34").Value = '  +2.79%  '

$ws.Range("D35").Value = "'0.0781"
$ws.Range("E35").Value = '  +2.34%  '

$ws.Range("E36").Value = '  +0.16%  '

$ws.Range("B37").Value = 'ARBITRUM'
$ws.Range("C37").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D37").Value = "'1.95"
$ws.Range("E37").Value = '  +3.60%  '

$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").Value = "'4.65"
$ws.Range("E38").Value = '  +4.52%  '

$ws.Range("E39").Value = '  +3.92%  '

$ws.Range("E40").Value = '  +1.83%  '

$ws.Range("D41").Value = "'121.27"
$ws.Range("E41").Value = '  -2.02%  '

$ws.Range("E42").Value = '  +0.72%  '

$ws.Range("D43").Value = "'21.34"
$ws.Range("E43").Value = '  +2.27%  '

$ws.Range("E44").Value = '  +2.23%  '

$ws.Range("D45").Value = "'1.969.65"
$ws.Range("E45").Value = '  +1.74%  '

$ws.Range("D46").Value = "'2.97"
$ws.Range("E46").Value = '  +1.33%  '

$ws.Range("E47").Value = '  -0.69%  '

$ws.Range("E48").Value = '  +2.90%  '

$ws.Range("D49").Value = "'9.22"
$ws.Range("E49").Value = '  -0.13%  '

$ws.Range("E50").Value = '  +13.90%  '

$ws.Range("D51").Value = "'78.86"
$ws.Range("E51").Value = '  +3.11%  '
